$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Lab 3 rpm1500.xlsx update: a new "TapPressure_H20_stddev" column is being
# inserted right after the existing TapPressure_H20 column (old column C,
# FlowRate_GPM, shifts right to column D), and a new "FlowRate_GPM_stddev"
# column is appended after it in column E.
# ---------------------------------------------------------------------------

# 1) Capture the existing FlowRate_GPM column (old column C) before it moves.
$oldC = $ws.Range("C1:C10").Value()

# 2) Write it, unchanged, into its new home at column D, then re-apply the
#    same centered alignment (style index 1) the original column carried.
$ws.Range("D1:D10").Value = $oldC
$ws.Range("D1:D10").HorizontalAlignment = -4108

# 3) Overwrite column C with the new TapPressure_H20_stddev header + data.
#    The header keeps the sheet's centered style, but the freshly authored
#    data values are left with no explicit style (unlike the other columns).
#    Column C used to hold FlowRate_GPM (centered, style index 1), so the
#    data cells need to be reset back to the plain "Normal" style first.
$ws.Range("C2:C10").Style = "Normal"
$ws.Range("C1").Value = "TapPressure_H20_stddev"
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C2").Value = 0.69934254839813692
$ws.Range("C3").Value = 0.29507626132916831
$ws.Range("C4").Value = 0.22711230702011945
$ws.Range("C5").Value = 0.13630847369110907
$ws.Range("C6").Value = 0.43637140144606268
$ws.Range("C7").Value = 0.20305171754998907
$ws.Range("C8").Value = 0.3246074552440219
$ws.Range("C9").Value = 0.14567086187703876
$ws.Range("C10").Value = 0.26884940022250448

# 4) Populate the brand-new FlowRate_GPM_stddev column E. Again the header
#    is centered like the rest of row 1, but the data cells are left plain.
#    (PowerShell-style "E-2" scientific literals aren't accepted by this
#    parser, so the same double values are written out in plain decimal.)
$ws.Range("E1").Value = "FlowRate_GPM_stddev"
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E2").Value = 0.035777087639997658
$ws.Range("E3").Value = 0.070213958726167924
$ws.Range("E4").Value = 0.074632432628181214
$ws.Range("E5").Value = 0.046151923036857709
$ws.Range("E6").Value = 0.049699094559157574
$ws.Range("E7").Value = 0.065954529791363861
$ws.Range("E8").Value = 0.060580524923443131
$ws.Range("E9").Value = 0.04472135954999603
$ws.Range("E10").Value = 0.077653074633268968

# 5) Match the widened/new column widths as closely as this engine's
#    pixel-quantized ColumnWidth property allows (it snaps to ~1/6-char
#    increments, so these land on the nearest achievable bucket).
$ws.Columns.Item(3).ColumnWidth = 21.0
$ws.Columns.Item(4).ColumnWidth = 13.83
$ws.Columns.Item(5).ColumnWidth = 19.25

# 6) Restore the last-selected cell reported after the edit.
$ws.Range("E22").Select()
